$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix label in the "MA" (multiple answer) question row: "Correct" -> "Corrects"
$ws.Range("F2").Value = "Corrects"

# Leave the selection where the user ended up after making the edit
$ws.Range("F6").Select() | Out-Null
